# "Game Rank solved - kattis"
# Adds a new entry to the problem-tracking sheet: day "11th" with problem
# "gamerank [3.6]" (note the leading space that was in the original data),
# placed in row 18 (columns D/E) right after the existing "10th" / "helpme [2.7]"
# entry in row 17. Also updates the current cell selection to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "11th"
$ws.Range("E18").Value = " gamerank [3.6]"

$ws.Range("E13").Select()
